# Report for Archive:
#  - Status text "Ready for handoff" -> "In Translation" on all three sheets
#    (Overview!E2/F2, zh-cn!C2, de-de!C2)
#  - Narrow the "status" column(s) that held that text:
#    Overview columns E & F, and column C on the zh-cn / de-de sheets,
#    from ~17.22 chars to ~13.41 chars.

$wb = $excel.ActiveWorkbook

# --- Update the status text wherever it appears ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# --- Narrow the corresponding columns ---
# (ColumnWidth of 12.5 "characters" is the value that lands the stored
#  column width nearest the ~13.41 target used for the archive report.)
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede.Columns.Item(3).ColumnWidth = 12.5
